$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.055.53'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.52%  '
$ws.Range("D3").Value = '''3.480.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.09%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''409.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("D6").Value = '''134.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +20.07%  '
$ws.Range("D7").Value = '''3.473.64'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.13%  '
$ws.Range("D8").Value = '''0.605'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.42%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '''0.697'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.80%  '
$ws.Range("D11").Value = '''0.132'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +33.60%  '
$ws.Range("D12").Value = '''43.60'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +10.34%  '
$ws.Range("D13").Value = '''0.143'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.67%  '
$ws.Range("D14").Value = '''4.030.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.80%  '
$ws.Range("D15").Value = '''8.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.23%  '
$ws.Range("D16").Value = '''20.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.71%  '
$ws.Range("D17").Value = '''3.466.44'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.39%  '
$ws.Range("D18").Value = '''63.019.73'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.67%  '
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").Value = '''10.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.58%  '
$ws.Range("D21").Value = '''0.0000140'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +27.04%  '
$ws.Range("D22").Value = '''3.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").Value = '''13.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.66%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '''82.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +9.79%  '
$ws.Range("D25").Value = '''315.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.57%  '
$ws.Range("D26").Value = '''3.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '''30.54'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.75%  '
$ws.Range("D28").Value = '''8.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.55%  '
$ws.Range("D29").Value = '''7.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("E30").Value = '  +0.52%  '
$ws.Range("D31").Value = '''4.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.18%  '
$ws.Range("E32").Value = '  +4.36%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").Value = '''44.07'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +11.65%  '
$ws.Range("E34").Value = '  +2.60%  '
$ws.Range("B35").Value = 'Cosmos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D35").Value = '''11.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.09%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = '''0.0495'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.21%  '
$ws.Range("D38").Value = '''52.61'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.56%  '
$ws.Range("D39").Value = '''3.62'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.29%  '
$ws.Range("E40").Value = '  -0.18%  '
$ws.Range("D41").Value = '''3.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.30%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = '''0.126'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.75%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '''2.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.14%  '
$ws.Range("D44").Value = '''137.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").Value = '''17.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.18%  '
$ws.Range("D46").Value = '''4.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.77%  '
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("E49").Value = '  +0.39%  '
$ws.Range("D50").Value = '''2.202.38'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '''3.825.06'
$ws.Range("D51").Style = "Normal"
